$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ship-to address block (E14:E17) — each line shifts up by one,
# and a new final line "REPLACE" is introduced.
$ws.Range("E14").Value = "UNFI " + [char]0x2013 + " Dublin, CA"
$ws.Range("E15").Value = "4000 Inspiration Drive"
$ws.Range("E16").Value = "Dublin, CA 94568"
$ws.Range("E17").Value = "REPLACE"

# Populate the previously empty production date field.
$ws.Range("E21").Value = "08/01/24"

# Update internal tracking / item numbers.
$ws.Range("C26").Value = "02420"
$ws.Range("C27").Value = "02421"
$ws.Range("C28").Value = "100262"
$ws.Range("C29").Value = "10068"
